$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the (now second-to-last-edit)
#    empty paragraph after "(A voir si ... )" to the end of the
#    paragraph that ends with "... identifier le séjour réservé par le
#    client." (i.e. where the author's last edit actually happened).
# ---------------------------------------------------------------------

$needle = "qui permettra d" + [char]0x2019 + "identifier le s" + [char]0xE9 + "jour r" + [char]0xE9 + "serv" + [char]0xE9 + " par le client."
$rng = $d.Content
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence for bookmark placement"
}

# $rng now covers the matched text; collapse to its end (just before the
# paragraph mark). Inserting a bookmark directly at that exact offset
# hits an engine edge case, so nudge past it with a throwaway character,
# bookmark there, then remove the throwaway character again.
$rng.Collapse(0)
$endPos = $rng.Start

$tmp = $d.Range($endPos, $endPos)
$tmp.InsertAfter("X")

$bmTarget = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmTarget) | Out-Null

$tmpChar = $d.Range($endPos, $endPos + 1)
$tmpChar.Delete()

# ---------------------------------------------------------------------
# 2) Remove the two empty paragraphs (one plain, one list-numbered)
#    that used to follow the "... statut du nettoyage." paragraph.
# ---------------------------------------------------------------------

$needle2 = "qui se mettra " + [char]0xE0 + " jour automatiquement avec le nombre de chambre qui auront le statut du nettoyage."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target sentence for paragraph cleanup"
}

$hostPara = $rng2.Paragraphs.First
$hostIndex = $hostPara.Index

# Delete the following two paragraphs one at a time (deleting them as a
# single combined range merges differently and under-deletes; re-fetch
# by index after each delete since the collection shifts).
$nextPara = $d.Paragraphs($hostIndex + 1)
$nextPara.Range.Delete()
$nextPara2 = $d.Paragraphs($hostIndex + 1)
$nextPara2.Range.Delete()
